# Apply updated forecast-error values to rows 2-10, and append new row 11 (Q9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B, C, D, E, F, G
$rows = @(
    @(2,  0.4257792306779735, 2.217877629047371, 14.65900008690008, 3.828707365012384, 3.842820107199375, 51),
    @(3,  0.7099789703697202, 2.408638386592436, 14.21756035022525, 3.770618032925802, 3.740769621219512, 50),
    @(4,  0.5489223891954144, 2.142174787177605, 13.77538404339546, 3.711520449006776, 3.708743454234583, 49),
    @(5,  0.7318547725335144, 2.680487997555033, 18.5702968988901,  4.309326733828627, 4.291666599746589, 48),
    @(6,  0.6680602970490725, 2.486664053058593, 14.69638859368018, 3.833586909629176, 3.815739570344253, 47),
    @(7,  0.6768152876208176, 2.661179530454853, 17.12325255943136, 4.138025200434545, 4.127409698292647, 46),
    @(8,  0.5558206604367482, 2.164092278013886, 13.2712682877262,  3.642975197242798, 3.641006746554708, 45),
    @(9,  0.6185392898960298, 2.459293877902776, 13.75218320165127, 3.70839361471396,  3.69871776441472,  44),
    @(10, 0.7423495001570053, 2.530096135466969, 17.10393720569436, 4.135690656431446, 4.116669771413056, 43)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
}

# New row 11 ("Q9"): copy formatting of the row-10 label cell onto A11, then fill values
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Q9"

$ws.Range("B11").Value = 0.5511858627957305
$ws.Range("C11").Value = 2.216421366742779
$ws.Range("D11").Value = 14.11284219713519
$ws.Range("E11").Value = 3.756706296363237
$ws.Range("F11").Value = 3.761095865066883
$ws.Range("G11").Value = 42
